$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.652167666666667
$ws.Range("H2").Value = 16.956503
$ws.Range("I2").Value = 0.1860329065948871
$ws.Range("J2").Value = 0.1860329065948871
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 87.94215800000001
$ws.Range("N2").Value = 263.826474
$ws.Range("O2").Value = 0.8507690866039653
$ws.Range("P2").Value = 0.8507690866039652
$ws.Range("Q2").Value = 497.0638219844914
$ws.Range("R2").Value = 4473.574397860422
$ws.Range("S2").Value = 0.1582710460220129
$ws.Range("T2").Value = 0.1582710460220128

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.652167666666667
$ws.Range("H3").Value = 16.956503
$ws.Range("I3").Value = 0.1860329065948871
$ws.Range("J3").Value = 0.1860329065948871
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 13.642319
$ws.Range("N3").Value = 40.926957
$ws.Range("O3").Value = 0.1319783769098539
$ws.Range("P3").Value = 0.1319783769098539
$ws.Range("Q3").Value = 77.10867435015234
$ws.Range("R3").Value = 693.9780691513711
$ws.Range("S3").Value = 0.02455232106421564
$ws.Range("T3").Value = 0.02455232106421564

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.652167666666667
$ws.Range("H4").Value = 16.956503
$ws.Range("I4").Value = 0.1860329065948871
$ws.Range("J4").Value = 0.1860329065948871
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.783357333333333
$ws.Range("N4").Value = 5.350072
$ws.Range("O4").Value = 0.01725253648618087
$ws.Range("P4").Value = 0.01725253648618087
$ws.Range("Q4").Value = 10.07983465757956
$ws.Range("R4").Value = 90.71851191821601
$ws.Range("S4").Value = 0.003209539508658567
$ws.Range("T4").Value = 0.003209539508658567

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 16.407289
$ws.Range("H5").Value = 49.221867
$ws.Range("I5").Value = 0.5400221369958743
$ws.Range("J5").Value = 0.5400221369958743
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 87.94215800000001
$ws.Range("N5").Value = 263.826474
$ws.Range("O5").Value = 0.8507690866039653
$ws.Range("P5").Value = 0.8507690866039652
$ws.Range("Q5").Value = 1442.892401589662
$ws.Range("R5").Value = 12986.03161430696
$ws.Range("S5").Value = 0.4594341402379014
$ws.Range("T5").Value = 0.4594341402379014

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 16.407289
$ws.Range("H6").Value = 49.221867
$ws.Range("I6").Value = 0.5400221369958743
$ws.Range("J6").Value = 0.5400221369958743
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 13.642319
$ws.Range("N6").Value = 40.926957
$ws.Range("O6").Value = 0.1319783769098539
$ws.Range("P6").Value = 0.1319783769098539
$ws.Range("Q6").Value = 223.833470463191
$ws.Range("R6").Value = 2014.501234168719
$ws.Range("S6").Value = 0.07127124513610623
$ws.Range("T6").Value = 0.07127124513610623

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 16.407289
$ws.Range("H7").Value = 49.221867
$ws.Range("I7").Value = 0.5400221369958743
$ws.Range("J7").Value = 0.5400221369958743
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.783357333333333
$ws.Range("N7").Value = 5.350072
$ws.Range("O7").Value = 0.01725253648618087
$ws.Range("P7").Value = 0.01725253648618087
$ws.Range("Q7").Value = 29.26005915826934
$ws.Range("R7").Value = 263.340532424424
$ws.Range("S7").Value = 0.009316751621866686
$ws.Range("T7").Value = 0.009316751621866686

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.323166333333333
$ws.Range("H8").Value = 24.969499
$ws.Range("I8").Value = 0.2739449564092387
$ws.Range("J8").Value = 0.2739449564092387
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 87.94215800000001
$ws.Range("N8").Value = 263.826474
$ws.Range("O8").Value = 0.8507690866039653
$ws.Range("P8").Value = 0.8507690866039652
$ws.Range("Q8").Value = 731.9572087462807
$ws.Range("R8").Value = 6587.614878716527
$ws.Range("S8").Value = 0.2330639003440511
$ws.Range("T8").Value = 0.2330639003440511

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.323166333333333
$ws.Range("H9").Value = 24.969499
$ws.Range("I9").Value = 0.2739449564092387
$ws.Range("J9").Value = 0.2739449564092387
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.642319
$ws.Range("N9").Value = 40.926957
$ws.Range("O9").Value = 0.1319783769098539
$ws.Range("P9").Value = 0.1319783769098539
$ws.Range("Q9").Value = 113.5472902093937
$ws.Range("R9").Value = 1021.925611884543
$ws.Range("S9").Value = 0.03615481070953199
$ws.Range("T9").Value = 0.036154810709532

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.323166333333333
$ws.Range("H10").Value = 24.969499
$ws.Range("I10").Value = 0.2739449564092387
$ws.Range("J10").Value = 0.2739449564092387
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.783357333333333
$ws.Range("N10").Value = 5.350072
$ws.Range("O10").Value = 0.01725253648618087
$ws.Range("P10").Value = 0.01725253648618087
$ws.Range("Q10").Value = 14.84317971710311
$ws.Range("R10").Value = 133.588617453928
$ws.Range("S10").Value = 0.004726245355655618
$ws.Range("T10").Value = 0.004726245355655619
